$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 140, shifting existing rows 140-200 down to 141-201.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new record.
$ws.Cells.Item(140, 1).Value = 6
$ws.Cells.Item(140, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(140, 3).Value = "Metropolitana"
$ws.Cells.Item(140, 4).Value = 44609
$ws.Cells.Item(140, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(140, 5).Value = 13
$ws.Cells.Item(140, 6).Value = 100112022
$ws.Cells.Item(140, 7).Value = "Arveja Verde"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 130
$ws.Cells.Item(140, 11).Value = 20000
$ws.Cells.Item(140, 12).Value = 23000
$ws.Cells.Item(140, 13).Value = 21154
$ws.Cells.Item(140, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(140, 15).Value = "Carahue"
$ws.Cells.Item(140, 16).Value = 846
$ws.Cells.Item(140, 17).Value = 25
$ws.Cells.Item(140, 18).Value = "Hortaliza"
